# Auto-generated edit script applying the diff changes to 上海-漫展信息.xlsx
# Commit message: Update gh-pages to output generated at 456a3b4

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (25 cell updates) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3270
$ws.Range("E5").Value = "2024.07.17 10:00-09.15 20:00"
$ws.Range("F8").Value = 7641
$ws.Range("F11").Value = 20
$ws.Range("F14").Value = 663
$ws.Range("F18").Value = 164
$ws.Range("F19").Value = 1553
$ws.Range("F20").Value = 330
$ws.Range("F21").Value = 6047
$ws.Range("F25").Value = 1003
$ws.Range("F27").Value = 4194
$ws.Range("F28").Value = 3840
$ws.Range("F29").Value = 287
$ws.Range("F30").Value = 90
$ws.Range("F31").Value = 1035
$ws.Range("F33").Value = 1024
$ws.Range("F34").Value = 1023
$ws.Range("F38").Value = 415
$ws.Range("F42").Value = 579
$ws.Range("F43").Value = 381
$ws.Range("F44").Value = 313
$ws.Range("F45").Value = 1089
$ws.Range("F46").Value = 464
$ws.Range("F47").Value = 3093
$ws.Range("F49").Value = 341

# --- Sheet: 演出 (13 cell updates) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = 604
$ws.Range("F12").Value = 76
$ws.Range("F15").Value = 241
$ws.Range("F20").Value = 150
$ws.Range("F21").Value = 27
$ws.Range("F22").Value = 33
$ws.Range("F26").Value = 24
$ws.Range("F28").Value = 5462
$ws.Range("F29").Value = 5462
$ws.Range("F30").Value = 44
$ws.Range("F33").Value = 52
$ws.Range("F35").Value = 1

# --- Sheet: 本地生活 (7 cell updates) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 580
$ws.Range("F6").Value = 1951
$ws.Range("F10").Value = 1262
$ws.Range("F12").Value = 521
$ws.Range("F13").Value = 2051
$ws.Range("F14").Value = 8753
$ws.Range("F16").Value = 55

# --- Sheet: 全部类型 (26 cell updates) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1951
$ws.Range("F7").Value = 7641
$ws.Range("F9").Value = 1262
$ws.Range("F10").Value = 521
$ws.Range("F13").Value = 20
$ws.Range("F17").Value = 604
$ws.Range("F18").Value = 55
$ws.Range("F19").Value = 663
$ws.Range("F23").Value = 76
$ws.Range("F24").Value = 164
$ws.Range("F25").Value = 241
$ws.Range("F26").Value = 330
$ws.Range("F27").Value = 6047
$ws.Range("F31").Value = 4194
$ws.Range("F32").Value = 287
$ws.Range("F33").Value = 1035
$ws.Range("F35").Value = 1024
$ws.Range("F36").Value = 1023
$ws.Range("F40").Value = 150
$ws.Range("F42").Value = 579
$ws.Range("F43").Value = 381
$ws.Range("F44").Value = 313
$ws.Range("F45").Value = 24
$ws.Range("F46").Value = 464
$ws.Range("F47").Value = 3093
$ws.Range("F49").Value = 5461
